# SIH-Divergents_DbSchema.xlsx update
# Author: Ruchi Pareek
# Updated database schema for tables user, application and applicationcomments
#
# Changes applied:
#  - User table: "userId - int auto increment" -> "Id - int auto increment" (B11)
#  - User table: "Username - email (varchar)" -> "userId - email (varchar)" (B14)
#  - User table: new row "userRole" (B16)
#  - Application table: "username -foreign key" -> "userId -foreign key" (D12)
#  - Application table: fill in blank row with "isActive - boolean" (D14)
#  - Application table: new row "dateOfSubmission" (D15)
#  - Application Comments table: new row "commentsId - varchar" (D25)
#  - Selection / scroll position of the sheet view updated

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Prepare formatting for the brand-new rows first (copy-format only, this
#     does not touch any cell values / the shared-string table) -------------
$ws.Range("D14").Copy()
$ws.Range("D15").PasteSpecial(-4122)

$ws.Range("B13").Copy()
$ws.Range("B16").PasteSpecial(-4122)

$ws.Range("F23").Copy()
$ws.Range("D25").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# --- Application table ------------------------------------------------------
# D14 already existed (blank, pre-formatted) - fill it in
$ws.Range("D14").Value = "isActive - boolean"

# D15 (new row)
$ws.Range("D15").Value = "dateOfSubmission"

# --- User table --------------------------------------------------------------
# B16 (new row)
$ws.Range("B16").Value = "userRole"

# B14: Username - email (varchar) -> userId - email (varchar)
$ws.Range("B14").Value = "userId - email (varchar)"

# B11: userId - int auto increment -> Id - int auto increment
$ws.Range("B11").Value = "Id - int auto increment"

# --- Application table (continued) -------------------------------------------
# D12: username -foreign key -> userId -foreign key
$ws.Range("D12").Value = "userId -foreign key"

# --- Application Comments table --------------------------------------------
# D25 (new row)
$ws.Range("D25").Value = "commentsId - varchar"

# --- Sheet view: selection ---------------------------------------------------
$ws.Range("D26").Select()
